$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal numbers need the column
# kept as literal text (matching the source data, which stores every
# Price/Volume cell as text) -- force text formatting before assignment.
$textCells = @("D4", "D5", "D7", "D10", "D11", "D13", "D17", "D20", "D22", "D25", "D29", "D30", "D31", "D36", "D38", "D41", "D43", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.307.50"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "1.647.52"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "217.49"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").Value = "20.03"
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("D11").Value = "0.0793"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.707.10"
$ws.Range("E12").Value = "  +4.03%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.31"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "1.874.98"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "63.56"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "26.283.29"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "196.78"
$ws.Range("E20").Value = "  +1.88%  "
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").Value = "10.07"
$ws.Range("E22").Value = "  +0.77%  "
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  -2.66%  "
$ws.Range("D25").Value = "143.17"
$ws.Range("E25").Value = "  +0.75%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  +1.96%  "
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").Value = "15.66"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").Value = "1.26"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").Value = "0.0508"
$ws.Range("E31").Value = "  +2.50%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  +1.75%  "
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("D36").Value = "0.916"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").Value = "1.142.61"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").Value = "0.557"
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("E40").Value = "  +0.61%  "
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("E42").Value = "  +2.64%  "
$ws.Range("D43").Value = "100.31"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").Value = "1.783.81"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("E46").Value = "  +1.71%  "
$ws.Range("D47").Value = "1.48"
$ws.Range("E47").Value = "  +2.52%  "
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("D49").Value = "7.74"
$ws.Range("E49").Value = "  +2.97%  "
$ws.Range("D50").Value = "0.417"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "0.0977"
$ws.Range("E51").Value = "  +2.15%  "
